# Update cryptos list price/volume figures per the source diff.
# Price ("D") cells whose new text parses as a plain number need the
# NumberFormat="@" / Style="Normal" dance so Excel keeps them as text
# (matching the original inlineStr cells) instead of silently coercing
# them into numeric cells and dropping trailing zeros / formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.321.22"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "3.717.13"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").Value = "3.713.57"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  +5.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "4.338.15"
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").Value = "3.715.88"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "68.197.14"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  +7.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "491.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.09%  "
$ws.Range("D34").Value = "3.860.00"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.108"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "3.662.78"
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.323"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "430.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("D50").Value = "2.761.36"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0351"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.82%  "